# Use Case Description(등록된 자전거 리스트 조회) - split the bike-detail-view
# step out of the delete-item use case so it reads as its own flow:
#   B4: "4. 상세내용 보기/항목 삭제 옵션을..." -> "4. 상세정보 보기/항목 삭제" 옵션...
#   A5: "5a. 상세내용 보기 버튼..."            -> "5a. 상세정보 보기" 버튼...
#   B5: "6a. ...상세내용(...)을 출력한다."      -> "6a. ...상세정보를 볼 수 있는 화면으로 이동한다."

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B4").Value = '4. "상세정보 보기/항목 삭제" 옵션을 출력한다.'
$ws.Range("A5").Value = '5a. "상세정보 보기" 버튼을 클릭한다.'
$ws.Range("B5").Value = '6a. 해당 자전거의 상세정보를 볼 수 있는 화면으로 이동한다.'

# Row 5 previously had an explicit (taller) custom height; the new, shorter
# text no longer needs it, so let Excel drop back to the default row height.
$ws.Rows.Item(5).AutoFit()

# Column B narrows slightly to fit the edited wording.
$ws.Columns.Item(2).ColumnWidth = 63

# Selection moves to B5 (the cell that was just edited last).
$ws.Range("B5").Select()

# Cosmetic: the saved tab ratio (sheet-tab/scrollbar split) moved from 560 to 550.
$excel.ActiveWindow.TabRatio = 550
